$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for first data row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-09 01:00:41"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-11-09 01:00:27"
$wsZhCn.Range("K2").Value = "2016-11-09 01:01:20"

# de-de sheet: "Correspond Handback DateTime" (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-11-09 01:01:38"
